# Update "want to go" counts (column F) for a set of conventions that appear
# in the "展览" (Exhibitions) sheet and the "全部类型" (All types) aggregate
# sheet. The same event can appear in more than one sheet, so updates are
# matched by the event name (column C) rather than by a fixed row number,
# which keeps the script correct even if row ordering differs slightly.

$wb = $excel.ActiveWorkbook

# Map of event name -> new value for column F ("想去人数")
$updates = @{
    "南昌·第二届漫拥动漫嘉年华mini" = 172
    "南昌·DSL国风动漫游戏嘉年华" = 174
    "南昌·New World国潮动漫博览会" = 4920
    "南昌·晨啼漫拥二次元随机舞蹈派对-热爱欢聚(免费活动)" = 25
    "赣州·十万伏特-次元音乐only" = 6
    "九江·第三届ACD动漫游戏嘉年华" = 538
    "吉安·COMIC LIFE次元假日04" = 500
    "景德镇·第十四届瓷都ACG动漫游戏博览会" = 1367
    "江西·广电·Unlimited Project 动漫游戏博览会" = 3395
    "江西·第二十二届九江ACJJ国际动漫展" = 398
    "赣州·COMIC WORLD次元创作同人季特典·SP·动漫游戏嘉年华" = 127
    "赣州·十万伏特-第六届青年文化综合展览会" = 112
    "抚州·临次元07国漫&运动番嘉年华" = 73
    "江西·ShiningStaR数字互娱嘉年华" = 2561
    "上饶·第一届星光次元国风动漫游戏嘉年华暨我和我的cos小伙伴们" = 125
    "新余·LD02国风动漫嘉年华" = 82
    "赣州·漫库书店次元漫展" = 40
    "江西·ShiningStaR数字互娱嘉年华 配音演员陈张太康、张惠霖专场见面会" = 124
    "南昌·Kpop New Life" = 41
}

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        # .Value2 is used for the read because .Value's getter does not
        # reliably round-trip string cell contents in this host; .Value2
        # (like Excel's real Value2) returns the raw string/number.
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
